$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.333.60"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.490.86"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.99"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.85"
$ws.Range("E6").Value = "  +7.57%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.478"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.69"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.089.01"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.79"
$ws.Range("E13").Value = "  +7.74%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.502.19"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.372.55"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.33"
$ws.Range("E19").Value = "  +6.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.48"
$ws.Range("E20").Value = "  +6.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.31"
$ws.Range("E21").Value = "  +3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.566"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.51"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +10.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.637.41"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.185"
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +10.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("E30").Value = "  +6.01%  "
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("E32").Value = "  +6.80%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.85"
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.62"
$ws.Range("E35").Value = "  +29.81%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").Value = "  +9.70%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.18"
$ws.Range("E37").Value = "  +5.51%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.80"
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  +10.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.529.85"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0769"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.802"
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("E43").Value = "  +8.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("E44").Value = "  +4.58%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +11.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.619.31"
$ws.Range("E47").Value = "  +6.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  +17.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.75"
$ws.Range("E49").Value = "  +7.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("E51").Value = "  +5.71%  "
